$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 43.193607330322266
$ws.Range("C2").Value = 6.603447914123535
$ws.Range("D2").Value = 27.661653518676758
$ws.Range("E2").Value = 57.85714340209961
